# Applies the cryptos-list refresh described in the commit:
#   "Updated cryptos list on Thu Jul 20 13:13:19 UTC 2023 with GitHub Actions"
#
# Column D ("Price") values that look numeric (e.g. "0.8116", "244.80") are
# written with a leading apostrophe so Excel stores them as literal text
# (quote-prefixed), matching the source workbook where every Price/Volume
# cell is a text string (note some prices use "." as a thousands separator,
# e.g. "30.266.99", which already round-trips as text with no extra work).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.266.99'
$ws.Range('E2').Value = '  +1.11%  '
# Row 3
$ws.Range('D3').Value = '1.918.54'
$ws.Range('E3').Value = '  +0.61%  '
# Row 4
$ws.Range('E4').Value = '  +0.04%  '
# Row 5
$ws.Range('D5').Value = '''0.8116'
$ws.Range('E5').Value = '  +1.96%  '
# Row 6
$ws.Range('D6').Value = '''244.80'
$ws.Range('E6').Value = '  +1.33%  '
# Row 7
$ws.Range('E7').Value = '  +0.06%  '
# Row 8
$ws.Range('D8').Value = '''0.3251'
$ws.Range('E8').Value = '  +3.12%  '
# Row 9
$ws.Range('D9').Value = '''27.00'
$ws.Range('E9').Value = '  +3.22%  '
# Row 10
$ws.Range('D10').Value = '''0.07225'
$ws.Range('E10').Value = '  +4.63%  '
# Row 11
$ws.Range('D11').Value = '''0.7900'
$ws.Range('E11').Value = '  +6.90%  '
# Row 12
$ws.Range('D12').Value = '''0.08095'
$ws.Range('E12').Value = '  +1.45%  '
# Row 13
$ws.Range('D13').Value = '1.918.11'
$ws.Range('E13').Value = '  +1.01%  '
# Row 14
$ws.Range('D14').Value = '''5.412'
$ws.Range('E14').Value = '  +4.37%  '
# Row 15
$ws.Range('D15').Value = '''93.88'
$ws.Range('E15').Value = '  +1.22%  '
# Row 16
$ws.Range('D16').Value = '30.286.13'
$ws.Range('E16').Value = '  +1.18%  '
# Row 17
$ws.Range('D17').Value = '''14.20'
$ws.Range('E17').Value = '  +1.93%  '
# Row 18
$ws.Range('D18').Value = '''6.063'
$ws.Range('E18').Value = '  +3.65%  '
# Row 19
$ws.Range('D19').Value = '''249.68'
$ws.Range('E19').Value = '  +2.14%  '
# Row 20
$ws.Range('D20').Value = '''0.000007849'
$ws.Range('E20').Value = '  +1.59%  '
# Row 21
$ws.Range('D21').Value = '2.176.37'
$ws.Range('E21').Value = '  +1.11%  '
# Row 22
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').Value = '''8.251'
$ws.Range('E22').Value = '  +21.36%  '
# Row 23
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '''1.001'
$ws.Range('E23').Value = '  +0.09%  '
# Row 24
$ws.Range('E24').Value = '  +0.04%  '
# Row 25
$ws.Range('D25').Value = '''0.1670'
$ws.Range('E25').Value = '  +19.18%  '
# Row 26
$ws.Range('D26').Value = '''9.481'
$ws.Range('E26').Value = '  +3.18%  '
# Row 27
$ws.Range('D27').Value = '''168.29'
$ws.Range('E27').Value = '  +0.34%  '
# Row 28
$ws.Range('D28').Value = '''19.00'
$ws.Range('E28').Value = '  +0.77%  '
# Row 29
$ws.Range('D29').Value = '''2.165'
$ws.Range('E29').Value = '  +7.04%  '
# Row 30
$ws.Range('D30').Value = '''1.386'
$ws.Range('E30').Value = '  +1.74%  '
# Row 31
$ws.Range('D31').Value = '''1.552'
$ws.Range('E31').Value = '  +2.37%  '
# Row 32
$ws.Range('E32').Value = '  +0.80%  '
# Row 33
$ws.Range('D33').Value = '''0.05836'
$ws.Range('E33').Value = '  +6.43%  '
# Row 34
$ws.Range('D34').Value = '''4.146'
$ws.Range('E34').Value = '  +1.67%  '
# Row 35
$ws.Range('D35').Value = '''1.294'
$ws.Range('E35').Value = '  +3.23%  '
# Row 36
$ws.Range('D36').Value = '''0.7480'
$ws.Range('E36').Value = '  +2.74%  '
# Row 37
$ws.Range('D37').Value = '''2.733'
$ws.Range('E37').Value = '  +0.38%  '
# Row 38
$ws.Range('D38').Value = '''0.9956'
$ws.Range('E38').Value = '  -0.37%  '
# Row 39
$ws.Range('D39').Value = '''0.01962'
$ws.Range('E39').Value = '  +2.26%  '
# Row 40
$ws.Range('D40').Value = '''2.818'
# Row 41
$ws.Range('D41').Value = '''0.4536'
$ws.Range('E41').Value = '  +2.88%  '
# Row 42
$ws.Range('D42').Value = '''74.65'
$ws.Range('E42').Value = '  +3.81%  '
# Row 43
$ws.Range('D43').Value = '''5.980'
$ws.Range('E43').Value = '  -2.51%  '
# Row 44
$ws.Range('D44').Value = '''0.8550'
$ws.Range('E44').Value = '  +2.29%  '
# Row 45
$ws.Range('E45').Value = '  +3.22%  '
# Row 46
$ws.Range('E46').Value = '  +0.08%  '
# Row 47
$ws.Range('D47').Value = '''103.39'
$ws.Range('E47').Value = '  +3.23%  '
# Row 48
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '''10.04'
# Row 49
$ws.Range('B49').Value = 'SynthetixNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D49').Value = '''3.113'
$ws.Range('E49').Value = '  +12.77%  '
# Row 50
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '1.017.98'
$ws.Range('E50').Value = '  +3.47%  '
# Row 51
$ws.Range('D51').Value = '''7.621'
$ws.Range('E51').Value = '  +1.61%  '
